$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B19").Value = "Ajout de l'affichage des listes des commandes et périphérique depuis la BDD"
$ws.Range("C19").NumberFormat = $ws.Range("C18").NumberFormat
$ws.Range("C19").Value = 42329
$ws.Range("D19").Value = 0.45

$ws.Range("F2:F6").Select()
